$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price + Volume(1h)) per the commit diff.
# D-column price cells are stored as literal TEXT in the source data (e.g. "20.555.62"
# or "0.9502" with significant trailing zeros), so force text via NumberFormat "@" before
# assigning, then reset the cell style back to Normal to avoid leaving a stray number format.

# Row 2: Bitcoin
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "20.555.62"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3: Ethereum
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.473.73"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.64%  "

# Row 5: USDC
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9502"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +5.87%  "

# Row 6: BNB
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "279.19"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7: XRP
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3644"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -1.11%  "

# Row 8: Cardano
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3057"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -2.45%  "

# Row 9: OKB
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "39.81"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +1.58%  "

# Row 10: Polygon
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "1.053"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +3.25%  "

# Row 11: Dogecoin
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06653"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +2.80%  "

# Row 12: BinanceUSD
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.006"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "

# Row 13: Polkadot
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "5.504"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "

# Row 14: Solana
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "17.99"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +3.37%  "

# Row 15: Chainlink
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "6.204"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

# Row 16: Dai
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9515"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +5.17%  "

# Row 17: ShibaInu
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001030"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18: WrappedEther
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "1.472.88"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19: TRON
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05937"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +5.96%  "

# Row 20: Litecoin
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "69.48"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +3.10%  "

# Row 21: Uniswap
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "5.478"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +0.37%  "

# Row 22: Avalanche
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "14.43"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "

# Row 23: Cosmos
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "11.06"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  -0.43%  "

# Row 25: WrappedBTC
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "20.606.63"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

# Row 26: Monero
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "143.12"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +5.86%  "

# Row 27: LidoDAOToken
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "2.113"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -3.56%  "

# Row 28: EthereumClassic
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "17.20"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "

# Row 29: WrappedliquidstakedEther2.0
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "1.633.42"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30: BitcoinCash
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "113.36"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "

# Row 31: HuobiToken
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "3.946"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +8.21%  "

# Row 32: Filecoin
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "5.006"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +2.47%  "

# Row 33: ImmutableX
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.8061"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34: Stellar
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07960"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "

# Row 35: WEMIXTOKEN
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.517"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +6.88%  "

# Row 36: TrustWalletToken
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "1.214"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +6.42%  "

# Row 37: Hedera
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05837"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -2.03%  "

# Row 38: InternetComputer(DFINITY)
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "4.715"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "

# Row 39: VeChain
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02052"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +2.80%  "

# Row 40: Aptos
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "10.33"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +1.02%  "

# Row 41: Frax
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9518"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +3.88%  "

# Row 42: Algorand
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1876"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +2.75%  "

# Row 43: FraxShare
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "7.376"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +10.83%  "

# Row 44: TheSandbox
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5295"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "

# Row 45: PancakeSwap
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "3.536"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46: EnergySwap
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "12.22"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "

# Row 47: Quant
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "117.84"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "

# Row 48: Decentraland
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5179"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "

# Row 49: NEARProtocol
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.810"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "

# Row 50: Cronos
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06459"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "

# Row 51: PaxDollar
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9820"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -1.87%  "
